# "cleaned up git cheat sheet"
#
# The cheat-sheet commands are lower-cased (Ls/Cd/Git/Def -> ls/cd/git/def),
# a new "git clone" line is inserted before "git init", the numbered
# Download/Upload steps get their leading "N) " pulled into its own run, and
# the _GoBack bookmark in the heading paragraph is moved to the front of the
# paragraph. Word's background proofer also re-segments runs around each
# word it flags (w:proofErr spellStart/spellEnd/gramStart/gramEnd).
#
# The simplest reliable way to reproduce that exact run/markup layout via
# COM is to hand the whole body as literal WordprocessingML to
# Range.InsertXML (the same "paste formatted XML" mechanism
# Word.Interop exposes), rather than trying to recreate every
# Find/Replace + proofing-pass step interactively.

$d = $word.ActiveDocument

$wordOpenXml = @'
<?xml version="1.0" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p><w:proofErr w:type="spellStart"/><w:r><w:t>l</w:t></w:r><w:r><w:t>s</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> – </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>lists</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> everything in directory</w:t></w:r></w:p><w:p><w:proofErr w:type="gramStart"/><w:r><w:t>c</w:t></w:r><w:r><w:t>d</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> – change of directory</w:t></w:r></w:p><w:p><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:t>g</w:t></w:r><w:r><w:t>it</w:t></w:r><w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> clone – creates new repository, connects it to central repository, and pulls in one command</w:t></w:r></w:p><w:p><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:t>g</w:t></w:r><w:r><w:t>it</w:t></w:r><w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>init</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> – creates new repository</w:t></w:r></w:p><w:p><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:t>d</w:t></w:r><w:r><w:t>ef</w:t></w:r><w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> (Repository) – Data base of every single version of uploaded files</w:t></w:r></w:p><w:p><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:t>g</w:t></w:r><w:r><w:t>it</w:t></w:r><w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> remote add origin – syncs with central online repository for pull/push ability for everyone</w:t></w:r></w:p><w:p><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:t>mkdir</w:t></w:r><w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> – creates new folder in repository</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Heading1"/></w:pPr><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:t xml:space="preserve">For File </w:t></w:r><w:r><w:t>Download/</w:t></w:r><w:r><w:t>Upload</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">1) </w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:t>g</w:t></w:r><w:r><w:t>it</w:t></w:r><w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> pull – downloads all changes that were made by all others participants, including new files.</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">2) </w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:t>g</w:t></w:r><w:r><w:t>it</w:t></w:r><w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> status – tells about current states of all files in Repo</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">3) </w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:t>git</w:t></w:r><w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> add – adds file/folders to NEXT (pending) COMMIT</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">4) </w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:t>git</w:t></w:r><w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> commit – m “ ’message here’ ” : commits to the LOCAL repository (Saving changes in GIT)</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">5) </w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:t>git</w:t></w:r><w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> push – uploads all LOCAL changes to the GLABAL Repo</w:t></w:r></w:p><w:p/><w:p/></w:body></w:document>
</pkg:xmlData></pkg:part></pkg:package>
'@

# Replacing the whole story range swaps in the new paragraphs verbatim,
# including the w:proofErr markers and the re-split runs, while leaving the
# section properties (page size/margins) at the end untouched.
$d.Content.InsertXML($wordOpenXml)
